$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Nodes_Example: add a new "Node Parent ID" column (D) capturing which node
# is the parent of each node (e.g. computers parented under their router).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Nodes_Example")

# Copy the formatting (border + alignment) used by the existing header/data
# columns onto the new column before filling in values, so the new cells
# pick up the same cell style already used across the sheet.
$ws1.Range("C1:C9").Copy()
$ws1.Range("D1:D9").PasteSpecial(-4122)

$ws1.Range("D1").Value = "Node Parent ID"

# Node 1 (Router A) and Node 4 (Router B) are top-level, so their parent id
# is left blank. Their children point back at them via Node ID.
$ws1.Range("D2").Value = $null
$ws1.Range("D3").Value = 1
$ws1.Range("D4").Value = 1
$ws1.Range("D5").Value = $null
$ws1.Range("D6").Value = 4
$ws1.Range("D7").Value = 4
$ws1.Range("D8").Value = 4
$ws1.Range("D9").Value = $null

$ws1.Columns.Item(4).ColumnWidth = 13.6

$ws1.Range("A1:D9").Select()
$ws1.Range("D9").Activate()

# ---------------------------------------------------------------------------
# Leftover selection changes on the "expected output" sheets from the same
# editing session.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Data Exchanges_Expected Output")
$ws5.Range("A1:F5").Select()

$ws6 = $wb.Worksheets.Item("Data Flows_Expected Output")
$ws6.Activate()
$ws6.Range("B1:F3").Select()
